# Auto-generated Excel COM-interop script applying numeric corrections
# to the per-leve profit/cost columns (H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 786.8333
$ws.Range("J2").Value = 1000.2857
$ws.Range("L2").Value = 1000.2857
$ws.Range("N2").Value = -1226.2857

$ws.Range("H9").Value = 7269.0625
$ws.Range("I9").Value = 10730.3
$ws.Range("J9").Value = 1500.3334
$ws.Range("K9").Value = 10730.3
$ws.Range("L9").Value = 1500.3334
$ws.Range("M9").Value = -10561.3
$ws.Range("N9").Value = -1838.3334

$ws.Range("H15").Value = 152098.34
$ws.Range("I15").Value = 152098.34
$ws.Range("K15").Value = 456295.02
$ws.Range("M15").Value = -456126.02

$ws.Range("H42").Value = 1939.3334
$ws.Range("I42").Value = 2305.2
$ws.Range("J42").Value = 110
$ws.Range("K42").Value = 6915.599999999999
$ws.Range("L42").Value = 330
$ws.Range("M42").Value = -6685.599999999999
$ws.Range("N42").Value = -790

$ws.Range("H80").Value = 878.94116
$ws.Range("I80").Value = 747.7
$ws.Range("J80").Value = 1066.4286
$ws.Range("K80").Value = 2243.1
$ws.Range("L80").Value = 3199.2858
$ws.Range("M80").Value = -1245.1
$ws.Range("N80").Value = -5195.2858

$ws.Range("H83").Value = 878.94116
$ws.Range("I83").Value = 747.7
$ws.Range("J83").Value = 1066.4286
$ws.Range("K83").Value = 6729.3
$ws.Range("L83").Value = 9597.857399999999
$ws.Range("M83").Value = -1737.3
$ws.Range("N83").Value = -19581.8574

$ws.Range("H111").Value = 508.125
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""

$ws.Range("H113").Value = 1702
$ws.Range("J113").Value = 1817.2
$ws.Range("L113").Value = 1817.2
$ws.Range("N113").Value = -8325.200000000001

$ws.Range("H125").Value = 31254674
$ws.Range("J125").Value = 31254674
$ws.Range("L125").Value = 281292066
$ws.Range("N125").Value = -281296986

$ws.Range("H131").Value = 1130.8334
$ws.Range("I131").Value = 997
$ws.Range("K131").Value = 2991
$ws.Range("M131").Value = 2049

$ws.Range("H132").Value = 2327.5
$ws.Range("I132").Value = 2398.1904
$ws.Range("K132").Value = 7194.5712
$ws.Range("M132").Value = -4664.5712

$ws.Range("H138").Value = 3791.848
$ws.Range("I138").Value = 2571.56
$ws.Range("J138").Value = 5244.5713
$ws.Range("K138").Value = 7714.68
$ws.Range("L138").Value = 15733.7139
$ws.Range("M138").Value = -2574.68
$ws.Range("N138").Value = -26013.7139

$ws.Range("H141").Value = 2421.6667
$ws.Range("I141").Value = 2102
$ws.Range("K141").Value = 6306
$ws.Range("M141").Value = -1126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2731.625
$ws.Range("J2").Value = 4601.706
$ws.Range("L2").Value = 4601.706
$ws.Range("N2").Value = -4827.706

$ws.Range("H45").Value = 1425.2
$ws.Range("I45").Value = 891.75
$ws.Range("J45").Value = 1780.8334
$ws.Range("K45").Value = 891.75
$ws.Range("L45").Value = 1780.8334
$ws.Range("M45").Value = -514.75
$ws.Range("N45").Value = -2534.8334

$ws.Range("H61").Value = 336341660
$ws.Range("I61").Value = 336341660
$ws.Range("K61").Value = 336341660
$ws.Range("M61").Value = -336341448

$ws.Range("H74").Value = 17860114
$ws.Range("I74").Value = 22730644
$ws.Range("K74").Value = 22730644
$ws.Range("M74").Value = -22729770

$ws.Range("H77").Value = 17860114
$ws.Range("I77").Value = 22730644
$ws.Range("K77").Value = 113653220
$ws.Range("M77").Value = -113648852

$ws.Range("H116").Value = 2731.625
$ws.Range("J116").Value = 4601.706
$ws.Range("L116").Value = 4601.706
$ws.Range("N116").Value = -9189.706

$ws.Range("H132").Value = 3035314.5
$ws.Range("I132").Value = 5005272
$ws.Range("J132").Value = 4610.385
$ws.Range("K132").Value = 15015816
$ws.Range("L132").Value = 13831.155
$ws.Range("M132").Value = -15013286
$ws.Range("N132").Value = -18891.155

$ws.Range("H136").Value = 336341660
$ws.Range("I136").Value = 336341660
$ws.Range("K136").Value = 1009024980
$ws.Range("M136").Value = -1009022430

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2731.625
$ws.Range("J3").Value = 4601.706
$ws.Range("L3").Value = 4601.706
$ws.Range("N3").Value = -4829.706

$ws.Range("H20").Value = 1049.6774
$ws.Range("I20").Value = 704.4286
$ws.Range("K20").Value = 704.4286
$ws.Range("M20").Value = -457.4286

$ws.Range("H22").Value = 4763352.5
$ws.Range("I22").Value = 568.2857
$ws.Range("K22").Value = 568.2857
$ws.Range("M22").Value = -395.2857

$ws.Range("H86").Value = 1358.381
$ws.Range("I86").Value = 1252
$ws.Range("J86").Value = 1996.6666
$ws.Range("K86").Value = 1252
$ws.Range("L86").Value = 1996.6666
$ws.Range("M86").Value = -129
$ws.Range("N86").Value = -4242.6666

$ws.Range("H89").Value = 1358.381
$ws.Range("I89").Value = 1252
$ws.Range("J89").Value = 1996.6666
$ws.Range("K89").Value = 6260
$ws.Range("L89").Value = 9983.333000000001
$ws.Range("M89").Value = -644
$ws.Range("N89").Value = -21215.333

$ws.Range("H134").Value = 22731506
$ws.Range("I134").Value = 29414890
$ws.Range("K134").Value = 88244670
$ws.Range("M134").Value = -88242135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3291.9604
$ws.Range("I122").Value = 3234.0322
$ws.Range("K122").Value = 9702.096600000001
$ws.Range("M122").Value = -7252.096600000001

$ws.Range("H134").Value = 41668080
$ws.Range("I134").Value = 55556828
$ws.Range("K134").Value = 166670484
$ws.Range("M134").Value = -166667949

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 174.44444
$ws.Range("I2").Value = 96
$ws.Range("J2").Value = 272.5
$ws.Range("K2").Value = 576
$ws.Range("L2").Value = 1635
$ws.Range("M2").Value = -463
$ws.Range("N2").Value = -1861

$ws.Range("H4").Value = 1335229
$ws.Range("I4").Value = 763294.4
$ws.Range("J4").Value = 3337000
$ws.Range("K4").Value = 2289883.2
$ws.Range("L4").Value = 10011000
$ws.Range("M4").Value = -2289771.2
$ws.Range("N4").Value = -10011224

$ws.Range("H39").Value = 2466.0833
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 3199.125
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 9597.375
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -10185.375

$ws.Range("H55").Value = 1937.25
$ws.Range("I55").Value = 125
$ws.Range("J55").Value = 3749.5
$ws.Range("K55").Value = 375
$ws.Range("L55").Value = 11248.5
$ws.Range("M55").Value = -198
$ws.Range("N55").Value = -11602.5

$ws.Range("H121").Value = 822738.9399999999
$ws.Range("I121").Value = 204199.6
$ws.Range("J121").Value = 1132008.6
$ws.Range("K121").Value = 612598.8
$ws.Range("L121").Value = 3396025.8
$ws.Range("M121").Value = -611288.8
$ws.Range("N121").Value = -3398645.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 419.7619
$ws.Range("I2").Value = 428.3158
$ws.Range("J2").Value = 338.5
$ws.Range("K2").Value = 428.3158
$ws.Range("L2").Value = 338.5
$ws.Range("M2").Value = -315.3158
$ws.Range("N2").Value = -564.5

$ws.Range("H63").Value = 38550
$ws.Range("I63").Value = 38550
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 38550
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -37864

$ws.Range("H66").Value = 38550
$ws.Range("I66").Value = 38550
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 115650
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -112218

$ws.Range("H132").Value = 25003402
$ws.Range("I132").Value = 41668000
$ws.Range("K132").Value = 125004000
$ws.Range("M132").Value = -125001470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3299.8
$ws.Range("I22").Value = 5550
$ws.Range("J22").Value = 1799.6666
$ws.Range("K22").Value = 5550
$ws.Range("L22").Value = 1799.6666
$ws.Range("M22").Value = -5255
$ws.Range("N22").Value = -2389.6666

$ws.Range("H27").Value = 3299.8
$ws.Range("I27").Value = 5550
$ws.Range("J27").Value = 1799.6666
$ws.Range("K27").Value = 5550
$ws.Range("L27").Value = 1799.6666
$ws.Range("M27").Value = -5443
$ws.Range("N27").Value = -2013.6666

$ws.Range("H46").Value = 1257
$ws.Range("I46").Value = 959.2
$ws.Range("K46").Value = 959.2
$ws.Range("M46").Value = -771.2

$ws.Range("H48").Value = 35791
$ws.Range("I48").Value = 28687
$ws.Range("K48").Value = 28687
$ws.Range("M48").Value = -28026

$ws.Range("H55").Value = 292.08334
$ws.Range("J55").Value = 382
$ws.Range("L55").Value = 382
$ws.Range("N55").Value = -728

